$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 56, shifting existing rows 56-102 down to 57-103
$ws.Rows.Item(56).Insert()

$ws.Cells.Item(56, 1).Value = 4
$ws.Cells.Item(56, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(56, 3).Value = 'Los Lagos'
$ws.Cells.Item(56, 4).Value = 44566
$ws.Cells.Item(56, 5).Value = 10
$ws.Cells.Item(56, 6).Value = 100112022
$ws.Cells.Item(56, 7).Value = 'Arveja Verde'
$ws.Cells.Item(56, 8).Value = 'Sin especificar'
$ws.Cells.Item(56, 9).Value = 'Primera'
$ws.Cells.Item(56, 10).Value = 35
$ws.Cells.Item(56, 11).Value = 27000
$ws.Cells.Item(56, 12).Value = 27000
$ws.Cells.Item(56, 13).Value = 27000
$ws.Cells.Item(56, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(56, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(56, 16).Value = 1080
$ws.Cells.Item(56, 17).Value = 25
$ws.Cells.Item(56, 18).Value = 'Hortaliza'
